$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reset header row formatting (was bold style s=1 with custom height) back to default
$ws.Rows.Item(1).ClearFormats()
$ws.Rows.Item(1).AutoFit()

# Write the full table of files/classes (A:B), rows 1-25
$ws.Cells.Item(1,1).Value = "file.pdf"
$ws.Cells.Item(1,2).Value = "COMPPROGRAMIN"
$ws.Cells.Item(2,1).Value = "קוח-שירות-לקוחות.pdf"
$ws.Cells.Item(2,2).Value = "OFFICEMANAGEMEN"
$ws.Cells.Item(3,1).Value = "קורות חיים מיכל לאער.pdf"
$ws.Cells.Item(3,2).Value = "COMPPROGRAMIN"
$ws.Cells.Item(4,1).Value = "קורות חיים אורי עוז מרזם.pdf"
$ws.Cells.Item(4,2).Value = "ARCHITECTUR"
$ws.Cells.Item(5,1).Value = "קורות חיים מעודכן 2021.pdf"
$ws.Cells.Item(5,2).Value = "ACCOUNTIN"
$ws.Cells.Item(6,1).Value = "קוח-לדוגמא-מזכירות (1).pdf"
$ws.Cells.Item(6,2).Value = "OFFICEMANAGEMEN"
$ws.Cells.Item(7,1).Value = "שירה טוויג.pdf"
$ws.Cells.Item(7,2).Value = "ACCOUNTIN"
$ws.Cells.Item(8,1).Value = "קורות חיים רחל רונס (1).pdf"
$ws.Cells.Item(8,2).Value = "ACCOUNTIN"
$ws.Cells.Item(9,1).Value = "קורות חיים  פניני הוקס 0534106650.pdf"
$ws.Cells.Item(9,2).Value = "GRAPHICSANDDESIG"
$ws.Cells.Item(10,1).Value = "CVשירה ג'יקובס  .pdf"
$ws.Cells.Item(10,2).Value = "COMPPROGRAMIN"
$ws.Cells.Item(11,1).Value = "אבוחצירא.pdf"
$ws.Cells.Item(11,2).Value = "COMPPROGRAMIN"
$ws.Cells.Item(12,1).Value = "קורות חיים ברוריה רבינו.pdf"
$ws.Cells.Item(12,2).Value = "COMPPROGRAMIN"
$ws.Cells.Item(13,1).Value = "שירה אלנקוה - הנדסאי תכנה .pdf"
$ws.Cells.Item(13,2).Value = "COMPPROGRAMIN"
$ws.Cells.Item(14,1).Value = "קורות חיים הדר המגיד 2021.pdf"
$ws.Cells.Item(14,2).Value = "COMPPROGRAMIN"
$ws.Cells.Item(15,1).Value = "קורות חיים לאה דנקור.pdf"
$ws.Cells.Item(15,2).Value = "COMPPROGRAMIN"
$ws.Cells.Item(16,1).Value = "קורות חיים לדוגמה למשרת הייטק.pdf"
$ws.Cells.Item(16,2).Value = "COMPPROGRAMIN"
$ws.Cells.Item(17,1).Value = "קורות חיים מרים ברודסקי.pdf"
$ws.Cells.Item(17,2).Value = "COMPPROGRAMIN"
$ws.Cells.Item(18,1).Value = "קורות חיים רחל בלעך.pdf"
$ws.Cells.Item(18,2).Value = "COMPPROGRAMIN"
$ws.Cells.Item(19,1).Value = "שירה בן תקוה הנדסאי תכנה (1) (1).pdf"
$ws.Cells.Item(19,2).Value = "COMPPROGRAMIN"
$ws.Cells.Item(20,1).Value = "אילה קוח חדש.pdf"
$ws.Cells.Item(20,2).Value = "EDUCATIO"
$ws.Cells.Item(21,1).Value = "קורות חיים רחל רונס (1).pdf"
$ws.Cells.Item(21,2).Value = "בתקשורת -- 0 --> תכנון -- 0 --> ב -- 0 --> לב -- 0 --> עיצוב -- 0 --> ACCOUNTIN"
$ws.Cells.Item(22,1).Value = "קורות חיים רחל רונס (1).pdf"
$ws.Cells.Item(22,2).Value = "בתקשורת -- 0 --> תכנון -- 0 --> ב -- 0 --> לב -- 0 --> עיצוב -- 0 --> ACCOUNTIN"
$ws.Cells.Item(23,1).Value = "אילה קוח חדש.pdf"
$ws.Cells.Item(23,2).Value = "בתקשורת -- 0 --> תכנון -- 0 --> ב -- 1 --> לידה -- 1 --> EDUCATIO"
$ws.Cells.Item(24,1).Value = "קורות חיים אבי.pdf"
$ws.Cells.Item(24,2).Value = "בתקשורת -- 0 --> תכנון -- 0 --> ב -- 0 --> לב -- 1 --> בס -- 0 --> ד -- 0 --> קורות -- 0 --> COMPPROGRAMING"
$ws.Cells.Item(25,1).Value = "קורות חיים אבי.pdf"
$ws.Cells.Item(25,2).Value = "בתקשורת -- 0 --> תכנון -- 0 --> ב -- 0 --> לב -- 1 --> בס -- 0 --> ד -- 0 --> קורות -- 0 --> COMPPROGRAMING"

# Highlight specific rows (class was newly / uncertainly assigned) in yellow
$ws.Range("A4:B4").Interior.Color = 65535
$ws.Range("A8:B8").Interior.Color = 65535
$ws.Range("A9:B9").Interior.Color = 65535
$ws.Range("A20:B20").Interior.Color = 65535

# Column widths: B widened slightly; C:F given widths even though unused
$ws.Columns.Item(2).ColumnWidth = 71.28571428571429
$ws.Columns.Item(3).ColumnWidth = 16.285714285714285
$ws.Columns.Item(4).ColumnWidth = 16.857142857142858
$ws.Columns.Item(5).ColumnWidth = 19.71428571428571
$ws.Columns.Item(6).ColumnWidth = 39.14285714285714

# Selection ends below the data, at A31
$ws.Range("A31").Select()
